$p = $ppt.ActivePresentation

# --- Slide 31 ("Code Challenge"): add a new bullet after
#     "Attempt to parse the weather data and load it into a case class."
#     and before the existing (empty, lvl=2) placeholder paragraph. ---
$s31 = $p.Slides.Item(31)
$shp31 = $s31.Shapes.Item(2)
$tr31 = $shp31.TextFrame.TextRange
$para3 = $tr31.Paragraphs(3)
$para3.InsertAfter([char]13 + "Try and build the rest of the application.")

# --- Slide 7 ("Motivations For This Evenings Discussion"): reorder the
#     three bullets so "Getting ready for the opportunities that IoT
#     presents." moves from last to first. ---
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(2)
$tr7 = $shp7.TextFrame.TextRange
$para1 = $tr7.Paragraphs(1)
$para2 = $tr7.Paragraphs(2)
$para3b = $tr7.Paragraphs(3)

# First move text through unique placeholders so the engine's run-diffing
# doesn't try to reuse overlapping text prefixes between paragraphs (which
# would otherwise split runs apart).
$para1.Text = "PLACEHOLDER_ONE_XXXXXXXXXXXXXXXXXXXXXXXXXXXXX"
$para2.Text = "PLACEHOLDER_TWO_YYYYYYYYYYYYYYYYYYYYYYYYYYYYY"
$para3b.Text = "PLACEHOLDER_THREE_ZZZZZZZZZZZZZZZZZZZZZZZZZZZ"

$para1.Text = "Getting ready for the opportunities that IoT presents."
$para2.Text = "Tired of working with Sandboxes"
$para3b.Text = "Tired of playing with human generated data"
